$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -6472.5345434786
$ws.Range("C2").Value = 13561.10628116392
$ws.Range("D2").Value = -6472.534545050638
$ws.Range("E2").Value = -616.0371926347913
$ws.Range("F2").Value = 55.96194096546753
$ws.Range("G2").Value = 45.63458162350963
$ws.Range("H2").Value = 56.24114139383937
$ws.Range("I2").Value = 45.42231789641141
$ws.Range("J2").Value = 55.9619409655719
$ws.Range("K2").Value = 45.63458162497881
$ws.Range("L2").Value = 54.24647791925514
$ws.Range("O2").Value = 45.62699102690868
$ws.Range("P2").Value = 54.24647792150091
$ws.Range("R2").Value = 8.953380639800754
$ws.Range("S2").Value = -17.90676128295929
$ws.Range("T2").Value = 8.953380643158548
$ws.Range("X2").Value = -92.33646748703023
$ws.Range("Y2").Value = -142.3360080380841
$ws.Range("Z2").Value = -92.33646748703015
$ws.Range("AE2").Value = -8.953380639800754
$ws.Range("AF2").Value = 8.953380643158548
$ws.Range("AG2").Value = 8.953380639800754
$ws.Range("AH2").Value = -17.90676128295929
$ws.Range("AI2").Value = 8.953380643158548
$ws.Range("AJ2").Value = 8.953380639800754
$ws.Range("AK2").Value = -8.953380643158548
$ws.Range("AL2").Value = 24.99977027552693
$ws.Range("AM2").Value = -24.99977027552692
$ws.Range("AN2").Value = -92.33646748703023
$ws.Range("AO2").Value = -142.3360080380841
$ws.Range("AP2").Value = -92.33646748703015
$ws.Range("AQ2").Value = -24.99977027552693
$ws.Range("AR2").Value = 24.99977027552692
$ws.Range("AS2").Value = 55.96194096546753
$ws.Range("AT2").Value = 55.96194096546753
$ws.Range("AU2").Value = 56.24114139383931
$ws.Range("AV2").Value = 56.24114139383931
$ws.Range("AW2").Value = 56.24114139383937
$ws.Range("AX2").Value = 55.9619409655719
$ws.Range("AY2").Value = 55.9619409655719
$ws.Range("AZ2").Value = 45.63458162350963
$ws.Range("BA2").Value = 45.63458162350963
$ws.Range("BB2").Value = 45.42231789641141
$ws.Range("BC2").Value = 45.4223178956417
$ws.Range("BD2").Value = 45.4223178971809
$ws.Range("BE2").Value = 45.63458162497881
$ws.Range("BF2").Value = 45.63458162497881
$ws.Range("BG2").Value = 54.24647791925514
$ws.Range("BJ2").Value = 45.62699102690868
$ws.Range("BK2").Value = 54.24647792150091